# [Fix] Projectile 오류 수정
# Rename the monster projectile prefab paths so they are namespaced under
# "Monster_" (fixes a projectile-path bug), and update the active selection
# on the SkillList sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SkillList")

# projectilePrefabPath (column B) values that used to read
# "/Projectiles/Monster/MeleeAttack" / "Fireball" / "Arrow" now need the
# "Monster_" prefix on the final path segment. Every row that referenced one
# of the three shared strings gets updated so the old strings are no longer
# used (they fall out of the shared-string table) and the new ones take
# their place.
$ws.Range("B2").Value  = "/Projectiles/Monster/Monster_MeleeAttack"
$ws.Range("B3").Value  = "/Projectiles/Monster/Monster_MeleeAttack"
$ws.Range("B4").Value  = "/Projectiles/Monster/Monster_MeleeAttack"
$ws.Range("B5").Value  = "/Projectiles/Monster/Monster_MeleeAttack"
$ws.Range("B6").Value  = "/Projectiles/Monster/Monster_MeleeAttack"
$ws.Range("B7").Value  = "/Projectiles/Monster/Monster_Fireball"
$ws.Range("B8").Value  = "/Projectiles/Monster/Monster_Arrow"
$ws.Range("B9").Value  = "/Projectiles/Monster/Monster_MeleeAttack"
$ws.Range("B10").Value = "/Projectiles/Monster/Monster_MeleeAttack"
$ws.Range("B11").Value = "/Projectiles/Monster/Monster_MeleeAttack"
$ws.Range("B12").Value = "/Projectiles/Monster/Monster_MeleeAttack"
$ws.Range("B13").Value = "/Projectiles/Monster/Monster_MeleeAttack"

# Move/save the active selection on the SkillList sheet to B9 (was F13).
$ws.Range("B9").Select() | Out-Null
